$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 16 new reconciliation rows (rows 86-101) with data that was
# auto-synced in this batch.
$newRows = @(
    ,@("237671823369", "MFS ENTREE COLLEGE MALANGUE", "Rte_1", "Cite Sic Stade Marion Ocm", 189785.8, 110633, -79152.79999999999, 0.5829361311541749, "Cite Sic")
    ,@("237672128028", "CAROLINE WAKO DJAMNOU", "Rte_1", "Cite Sic Stade Marion Ocm", 25030, 22483, -2547, 0.8982421094686376, "Cite Sic")
    ,@("237672277367", "TOP MOBIL KM5 LTDLA_POLAS_BTQ_KM5", "Rte_0", "Cite Sic Stade Marion Ocm", 26250, 726946, 700696, 27.69318095238095, "Cite Sic")
    ,@("237674853971", "NJOSSEU TCHOUNZOU TOP MOBILE", "Rte_1", "Cite Sic Stade Marion Ocm", 182169.7, 129189, -52980.70000000001, 0.7091684292173726, "Cite Sic")
    ,@("237674884705", "manuela verna yetna baaga", "Rte_0", "Cite Sic Stade Marion Ocm", 5000, 15105, 10105, 3.021, "Cite Sic")
    ,@("237675779272", "RODES NGWEM KEMAYOU", "Rte_1", "Cite Sic Stade Marion Ocm", 184162.0370370371, 28759, -155403.0370370371, 0.1561613916890822, "Cite Sic")
    ,@("237677304210", "FERDINAND NKWELLE NGOME", "Rte_2", "Cite Sic Stade Marion Ocm", 61193.33333333334, 137078, 75884.66666666666, 2.240080618803791, "Cite Sic")
    ,@("237678267353", "LA NEGRESSE SARL EMBOLA BELTUS MBU", "Rte_0", "Cite Sic Stade Marion Ocm", 16520, 0, -16520, 0, "Cite Sic")
    ,@("237678370615", "ESSEN ONGOLONG BERTHE HORTENSE ETS MOBILE FINANCIAL SERVICES MFS", "Rte_1", "Cite Sic Stade Marion Ocm", 100229.8, 239366, 139136.2, 2.388171980788149, "Cite Sic")
    ,@("237678836319", "KAMDOM DOMINIQUE STEPHANIE ETS MOBILE FINANCIAL SERVICES MFS", "Rte_1", "Cite Sic Stade Marion Ocm", 25619.09090909091, 94879, 69259.90909090909, 3.703449132394166, "Cite Sic")
    ,@("237678922502", "NWOAGA TCHAMDJOU EPSE KAMSEU EMILINE ETS LE CONTENT", "Rte_0", "Cite Sic Stade Marion Ocm", 56220, 1196062, 1139842, 21.2746709356101, "Cite Sic")
    ,@("237679884264", "MFS CICAM", "Rte_1", "Cite Sic Stade Marion Ocm", 59610, 66267, 6657, 1.111675893306492, "Cite Sic")
    ,@("237681019523", "ETS MOULAY RIPERT AND COMPANY", "Rte_1", "Cite Sic Stade Marion Ocm", 57805.71428571429, 132723, 74917.28571428571, 2.296018683274021, "Cite Sic")
    ,@("237681125655", "EMENGUE PICHOU ROMEO KAMILAH CONNECTION GROUP", "Rte_1", "Cite Sic Stade Marion Ocm", 47813.75, 498608, 450794.25, 10.42812998352984, "Cite Sic")
    ,@("237681240793", "MBANE EMILIE FRANCOISE ETS MOBILE FINANCIAL SERVICES MFS", "Rte_0", "Cite Sic Stade Marion Ocm", 7438.090909090909, 15143, 7704.909090909091, 2.035871863503587, "Cite Sic")
    ,@("237682117915", "MEKUEKO FOUDJO BERLINE DIDIANE ETS MOBILE FINANCIAL SERVICES MFS", "Rte_1", "Cite Sic Stade Marion Ocm", 123740, 272599, 148859, 2.202998222078552, "Cite Sic")
)

$startRow = $ws.UsedRange.Rows.Count + 1
for ($idx = 0; $idx -lt $newRows.Count; $idx++) {
    $r = $startRow + $idx
    $data = $newRows[$idx]
    $ws.Cells.Item($r, 1).Value = [double]$data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = [double]$data[4]
    $ws.Cells.Item($r, 6).Value = [double]$data[5]
    $ws.Cells.Item($r, 7).Value = [double]$data[6]
    $ws.Cells.Item($r, 8).Value = [double]$data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
}

Write-Output "Added $($newRows.Count) rows"
